$d = $word.ActiveDocument

$find = "Kampagnendaten Perseus-Konstellation 2022: 16.-25. Januar, 7.-16. November, 6.-15. Dezember"
$replace = "Kampagnendaten 2022 für das Sternbild Perseus-Konstellation: 16.-25. Januar, 7.-16. November, 6.-15. Dezember"

$range = $d.Content
$range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
